# "Limitation continued" slide: merge the two runs of the second bullet
# ("still expects " + "8 electives") into a single run and drop the
# now-stray <a:endParaRPr/>, matching the author's edit.
#
# The host's TextRange.Text setter only rewrites the minimal span that
# differs from the current text (it keeps a common prefix/suffix run
# untouched), so a direct "merge these two runs" assignment would just
# leave the original run split in place. To force a real, from-scratch
# rewrite of the whole text body (which coalesces same-format runs into
# one new run carrying the first run's rPr, and drops stray endParaRPr),
# we first overwrite the whole text with unrelated placeholder content,
# then set it back to the real desired text in a second pass.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$line1 = "Does not currently deal with electives that are 3 credits"
$line2 = "still expects 8 electives"

# Pass 1: blow away the existing run/paragraph structure with placeholder
# text that shares no prefix/suffix with the real content, forcing the
# host to fully rebuild the text body instead of patching around it.
$tr.Text = "QQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQQ"

# Pass 2: write the real two-line content. Both paragraphs are rebuilt as
# single coalesced runs, and the previously-present endParaRPr on the
# second paragraph is gone.
$tr2 = $shp.TextFrame.TextRange
$tr2.Text = $line1 + "`r" + $line2

# The rebuild above drops the outline level on the second paragraph, so
# restore it (lvl="1" in OOXML <-> IndentLevel 2 in the COM model).
$para2 = $tr2.Paragraphs(2, 1)
$para2.IndentLevel = 2
